$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row (row 16) of data mirroring the existing rows' structure.
$row = 16

$ws.Range("A" + ($row - 1)).Copy()
$ws.Range("A" + $row).PasteSpecial(-4122)  # xlPasteFormats
$ws.Cells.Item($row, 1).Value = 42622.887291666666

$ws.Cells.Item($row, 2).Value = 58
$ws.Cells.Item($row, 3).Value = 0
$ws.Cells.Item($row, 4).Value = 0
$ws.Cells.Item($row, 5).Value = 0
$ws.Cells.Item($row, 6).Value = 0
$ws.Cells.Item($row, 7).Value = 0
$ws.Cells.Item($row, 8).Value = 0
$ws.Cells.Item($row, 9).Value = 0
$ws.Cells.Item($row, 10).Value = 0
$ws.Cells.Item($row, 11).Value = 0
$ws.Cells.Item($row, 12).Value = 0
$ws.Cells.Item($row, 13).Value = 0
$ws.Cells.Item($row, 14).Value = "Random"
